# expansão das análises automáticas
# Adds three new computed-metric columns (L: apoio_medio, M: contribuicoes,
# N: media_contribuicoes) to the summary table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the formatting of the existing header cells (bold, centered,
# bordered) by copying the format from K1 onto the new header cells.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-7) ------------------------------------------------
$values = @{
    2 = @(91.11272172566387, 192790, 304.565560821485)
    3 = @(91.82281991764464, 70763,  359.2030456852792)
    4 = @(89.11804758541258, 143785, 129.536036036036)
    5 = @(92.68628097576973, 59861,  219.2710622710623)
    6 = @(18.00852858651895, 1674,   12.77862595419847)
    7 = @(24.42075075084659, 534,    25.42857142857143)
}

foreach ($row in $values.Keys) {
    $triple = $values[$row]
    $ws.Cells.Item($row, 12).Value = $triple[0]  # L -> apoio_medio
    $ws.Cells.Item($row, 13).Value = $triple[1]  # M -> contribuicoes
    $ws.Cells.Item($row, 14).Value = $triple[2]  # N -> media_contribuicoes
}
